# Finding Anagrams v1.0 - update the title-slide date from
# "February 20, 2021" to "February 22, 2021".
#
# The author's edit (per the OOXML diff) retyped the day-of-month, which
# caused PowerPoint to split the single date run into two runs:
#   1) "February 22, "
#   2) "2021"
# We reproduce that by editing the run text in place via a Characters()
# sub-range, which splits the run at the same boundary.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the subtitle shape (contains the author/date block) on slide 1.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -like "*February 20, 2021*") {
            $shape = $candidate
        }
    }
}

$tr = $shape.TextFrame.TextRange

# Find the paragraph holding the date text.
$datePara = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $candidatePara = $tr.Paragraphs($i)
    if ($candidatePara.Text -eq "February 20, 2021") {
        $datePara = $candidatePara
    }
}

# Replace just the "February 20, " portion (leaving the year "2021"
# untouched) so the edit splits into two runs at the same place a live
# retype of the day would: "February 22, " | "2021".
$prefixLen = "February 20, ".Length
$prefix = $datePara.Characters(1, $prefixLen)
$prefix.Text = "February 22, "
